$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = @'
questions = [
    {
        "title": "You have been outsourced by a prestigious university and need to analyze the quality of a new online course. \u200b\u200bTo carry out the research, the program manager asks you to survey only students who received high grades in the class, which is a minority.  What recommendation should you make to the program manager?",
        "ques_type": 2,
        "options": [
            "The survey should be applied to the same number of students with low grades as those with high grades.",
            "The survey should be applied only to the students who received a low grade to determine if the problem is related to the course.",
            "The survey should be applied to the course coordinators since they have the most direct, in-depth experience with the course.",
            "The survey should be applied to all participating students, including feedback from the coordinators."
        ],
        "score": "The survey should be applied to all participating students, including feedback from the coordinators."
    },
    {
        "title": "You are the market researcher designing a questionnaire for understanding customer pain points for a software as a service (SaaS) company. The marketing manager wants to know how effective the company\u2019s live chat support is at troubleshooting problems. The below chart contains a list of questions you are considering using for the survey. Which option should you recommend the marketing manager use?",
        "ques_type": 2,
        "options": [
            "Option A",
            "Option B",
            "Option C",
            "Option D"
        ],
        "score": "Option C"
    },
    {
        "title": "You were outsourced by a consumer packaged goods company to help with a pre-launch qualitative research interview for a new kitchen cleaning product. The primary stakeholders asked you to determine whether this new product is going to be viable.  Which of the provided questions should you include in the interview?",
        "ques_type": 2,
        "options": [
            "What do you like most about cleaning your kitchen?",
            "Please describe the most important factors about cleaning your home.",
            "Would you buy this cleaning product?",
            "Why do you always want to keep your kitchen clean?"
        ],
        "score": "What do you like most about cleaning your kitchen?"
    },
    {
        "title": "You are the market research manager of a consulting firm. One of your clients requested an anonymous survey to be applied to their company\u2019s senior marketing managers to find out what services they might be interested in from an external multimedia agency. The client who requested this research asks you for the personal information of people who chose a specific answer.  True or false: The best course of action is to charge an additional fee for this information.",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "False"
    }
]
'@

$ws.Range("A2").Clear()
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $questionsText
